# Insert a new "Industry" column at C (shifting Mutual Fund..QoQ from C:I to D:J)
# and populate it with the industry classification for each holding row.
# Source: motilal_portfolio_change_engine latest processed run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C. Excel shifts the existing C:I data
# (Mutual Fund, Status, Jan_2026, Dec_2025, Oct_2025, MoM, QoQ) to D:J and
# carries the header style (bold/centered/bordered) along for the new cell.
$ws.Columns.Item(3).EntireColumn.Insert()

$ws.Range("C1").Value = "Industry"

$industries = @{
    2  = "Healthcare Services"
    3  = "Transport Services"
    4  = "Retailing"
    5  = "Capital Markets"
    6  = "Finance"
    7  = "Finance"
    8  = "Retailing"
    9  = "Banks"
    10 = "Pharmaceuticals & Biotechnology"
    11 = "Banks"
    12 = "Retailing"
    13 = "Banks"
    14 = "Realty"
    15 = "Pharmaceuticals & Biotechnology"
    16 = "Capital Markets"
    17 = "IT - Software"
    18 = "Cement & Cement Products"
    19 = "Consumer Durables"
    20 = "Electrical Equipment"
    21 = "Healthcare Services"
    22 = "Insurance"
    23 = "Healthcare Services"
    24 = "Industrial Manufacturing"
    25 = "Automobiles"
    26 = "Healthcare Services"
    27 = "Cement & Cement Products"
    28 = "Realty"
    29 = "Consumer Durables"
    30 = "Industrial Products"
    31 = "Consumer Durables"
    32 = "IT - Software"
    33 = "Insurance"
    34 = "Industrial Products"
    35 = "Realty"
    36 = "Realty"
    37 = "Auto Components"
    38 = "Industrial Products"
    39 = "Industrial Products"
    40 = "Consumer Durables"
    41 = "Commercial Services & Supplies"
    42 = "Pharmaceuticals & Biotechnology"
    43 = "Retailing"
    44 = "Automobiles"
    45 = "Insurance"
    46 = "Banks"
}

foreach ($row in 2..46) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
